# Update weight tracking sheet with new data (images of the scale) for
# weeks 36-40, fix week 7 value, and extend the table with 4 new (empty)
# weeks 41-44.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("weight")   # the "weight" tab (already the active sheet)

# --- Fix existing entry: week 7 (row 10) 52 -> 51 ---
$ws.Range("E10").Value = 51

# --- Row 39 (week 36): fill in Thu/Fri/Sat/Sun (H:K) ---
$ws.Range("H39").Value = 60.1
$ws.Range("I39").Value = 60.1
$ws.Range("J39").Value = 60.3
$ws.Range("K39").Value = 60.3

# --- Row 40 (week 37): fill in the whole week ---
$ws.Range("E40").Value = 60.1
$ws.Range("F40").Value = 60.4
$ws.Range("G40").Value = 60.5
$ws.Range("H40").Value = 60.7
$ws.Range("I40").Value = 60.8
$ws.Range("J40").Value = 60.2
$ws.Range("K40").Value = 60.1

# --- Row 41 (week 38): fill in the whole week ---
$ws.Range("E41").Value = 60.3
$ws.Range("F41").Value = 60.9
$ws.Range("G41").Value = 60.9
$ws.Range("H41").Value = 60.9
$ws.Range("I41").Value = 60.8
$ws.Range("J41").Value = 60.7
$ws.Range("K41").Value = 60.5

# --- Row 42 (week 39): only first day recorded ---
$ws.Range("E42").Value = 61

# --- Extend the table with 4 more weeks (rows 44-47), column A only ---
$ws.Range("A43").Copy($ws.Range("A44"))
$ws.Range("A43").Copy($ws.Range("A45"))
$ws.Range("A43").Copy($ws.Range("A46"))
$ws.Range("A43").Copy($ws.Range("A47"))
$ws.Range("A44:A47").Formula = "=A43+7"

# --- Update the frozen pane / view so it now shows the top of the table ---
$ws.Activate()
$ws.Range("B3").Select()
$ws.Range("H10").Select()

$wb.Save()
